$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Complete row 4 (week 3) with the actual/prediction data and derived formulas
$ws.Range("D4").Value = 5.9861603567514203
$ws.Range("E4").Value = 6.3418596432486396
$ws.Range("F4").Value = 6.1640100000000002
$ws.Range("G4").Value = 6.58
$ws.Range("H4").Formula = "=G4-F4"
$ws.Range("I4").Formula = "=H4/G4"

# Add row 5 (week 4) following the same pattern as the previous weeks
$ws.Range("A5").Value = 4
$ws.Range("B4:C4").Copy($ws.Range("B5:C5"))
$ws.Range("B5").Formula = "=B4+7"
$ws.Range("C5").Formula = "=C4+7"
$ws.Range("D5").Value = 5.9528503567514202
$ws.Range("E5").Value = 6.3085496432486403
$ws.Range("F5").Value = 6.1307

$ws.Range("G5").Select()
